$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SchemeMaster")

# ---------------------------------------------------------------
# 1. Insert two new columns at the front (A, B) for TC_ID / TC_Name.
#    Everything that used to be in columns A..N shifts right to C..P.
# ---------------------------------------------------------------
$ws.Columns("A:B").Insert()

# ---------------------------------------------------------------
# 2. New header cells for the inserted columns (copy style from the
#    neighbouring header cell first, then overwrite the text).
# ---------------------------------------------------------------
$ws.Range("C1").Copy($ws.Range("A1"))
$ws.Range("C1").Copy($ws.Range("B1"))
$ws.Range("A1").Value2 = "TC_ID"
$ws.Range("B1").Value2 = "TC_Name"

# ---------------------------------------------------------------
# 3. Fill in the Test-Case columns for the existing data row (row 2)
#    and duplicate the whole row into a brand-new row 3 for the
#    second test case, then customise the bits that differ.
# ---------------------------------------------------------------
$ws.Range("C2").Copy($ws.Range("A2"))
$ws.Range("C2").Copy($ws.Range("B2"))
$ws.Range("A2").Value2 = "TC_01"
$ws.Range("B2").Value2 = "Validate Duplicate Scheme Creation"

# Duplicate row 2 (now fully populated C..P) down into row 3, carrying
# across both values and formatting.
$ws.Range("A2:P2").Copy($ws.Range("A3"))

$ws.Range("A3").Value2 = "TC_02"
$ws.Range("B3").Value2 = "Validate New Scheme Creation"
$ws.Range("C3").Value2 = "Scheme Test"
$ws.Range("E3").Value2 = "This is scheme description for "

# ---------------------------------------------------------------
# 4. Column widths.
#    Columns that merely shifted keep their original width automatically;
#    only the genuinely new/changed ones need to be set explicitly.
#    (engine stores width = ColumnWidth + 5/6, rounded to the nearest 1/6)
# ---------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 19.706666666666667   # -> 20.54
$ws.Columns("B").ColumnWidth = 30.28666666666667    # -> 31.12
$ws.Columns("E").ColumnWidth = 25.886666666666667   # -> 26.72

# ---------------------------------------------------------------
# 5. Data validations now need to cover both data rows (2 and 3).
#    Column insertion already shifted the sqref letters; extend them
#    to span row 3 as well by recreating the rules.
# ---------------------------------------------------------------
$ws.Cells.Validation.Delete()

$v1 = $ws.Range("D2:D3").Validation
$v1.Add(3, 1, 1, """Capri Gold Loans,Shivalik Small Finanace Bank,AGRI,MSME,Bank Of Baroda,Karur Vysya Bank""")
$v1.IgnoreBlank = $true
$v1.InCellDropdown = $true
$v1.ShowInput = $false
$v1.ShowError = $true

$v2 = $ws.Range("F2:F3").Validation
$v2.Add(3, 1, 1, """Standard,Rebate""")
$v2.IgnoreBlank = $true
$v2.InCellDropdown = $true
$v2.ShowInput = $false
$v2.ShowError = $true

$v3 = $ws.Range("M2:M3").Validation
$v3.Add(3, 1, 1, """Monthly,Bi-Monthly,Quarterly,Half Yearly,Yearly""")
$v3.IgnoreBlank = $true
$v3.InCellDropdown = $true
$v3.ShowInput = $false
$v3.ShowError = $true

$v4 = $ws.Range("N2:N3").Validation
$v4.Add(3, 1, 1, """MTM Charges,Processing Fee,SOA Charges,Pre-Auction Charges,Post-Auction Charges,Courier Charges,Other Charges,Legal Charges""")
$v4.IgnoreBlank = $true
$v4.InCellDropdown = $true
$v4.ShowInput = $false
$v4.ShowError = $true

$v5 = $ws.Range("O2:O3").Validation
$v5.Add(3, 1, 1, """Flat Value,Percentage Of Loan Amount""")
$v5.IgnoreBlank = $true
$v5.InCellDropdown = $true
$v5.ShowInput = $false
$v5.ShowError = $true

# ---------------------------------------------------------------
# 6. Match the selection recorded in the saved file.
# ---------------------------------------------------------------
$ws.Range("E4").Select()
